$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65
$ws.Cells.Item($row, 1).Value = "Debiasi Alessio"
$ws.Cells.Item($row, 2).Value = "Elia Battisti | U.SGUARNA"
$ws.Cells.Item($row, 3).Value = "Michele Merighi | Clitoriders"
$ws.Cells.Item($row, 4).Value = "Amedeo Malesardi | FC SAVIGNANO"
$ws.Cells.Item($row, 5).Value = "Federico Nicolodi | U.SGUARNA"
$ws.Cells.Item($row, 6).Value = "Alessio Debiasi | Mai una gioia"
